# Updated cryptos list with latest prices and percentage changes.
# Cells in column D that are plain numbers must be forced to Text format
# first, since the source data in this sheet is stored as text (inline
# strings) rather than numeric values -- assigning a plain numeric-looking
# string to .Value would otherwise be auto-coerced into a Double by Excel.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = '27.006.24'
$ws.Range("E2").Value = '  +0.85%  '
# Row 3
$ws.Range("D3").Value = '1.557.99'
$ws.Range("E3").Value = '  +0.63%  '
# Row 4
$ws.Range("E4").Value = '  +0.57%  '
# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '207.53'
$ws.Range("E5").Value = '  +0.53%  '
# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.488'
$ws.Range("E6").Value = '  +1.41%  '
# Row 7
$ws.Range("E7").Value = '  +0.51%  '
# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '21.71'
$ws.Range("E8").Value = '  +1.43%  '
# Row 9
$ws.Range("E9").Value = '  +1.35%  '
# Row 10
$ws.Range("E10").Value = '  +1.61%  '
# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0863'
$ws.Range("E11").Value = '  +1.05%  '
# Row 12
$ws.Range("D12").Value = '1.779.14'
$ws.Range("E12").Value = '  +0.61%  '
# Row 13
$ws.Range("D13").Value = '1.557.28'
$ws.Range("E13").Value = '  +0.73%  '
# Row 15
$ws.Range("E15").Value = '  +1.05%  '
# Row 16
$ws.Range("B16").Value = 'Litecoin'
$ws.Range("C16").Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '62.06'
$ws.Range("E16").Value = '  +1.56%  '
# Row 17
$ws.Range("B17").Value = 'WrappedBTC'
$ws.Range("C17").Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
$ws.Range("D17").Value = '27.009.16'
$ws.Range("E17").Value = '  +0.90%  '
# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '216.18'
# Row 19
$ws.Range("D19").Value = '0.0₃0690'
$ws.Range("E19").Value = '  +0.59%  '
# Row 20
$ws.Range("E20").Value = '  +0.56%  '
# Row 21
$ws.Range("E21").Value = '  +0.45%  '
# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '4.04'
$ws.Range("E22").Value = '  -0.93%  '
# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '9.23'
$ws.Range("E23").Value = '  +3.05%  '
# Row 24
$ws.Range("E24").Value = '  -0.95%  '
# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '152.50'
$ws.Range("E25").Value = '  -0.45%  '
# Row 26
$ws.Range("E26").Value = '  +2.29%  '
# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '14.94'
$ws.Range("E27").Value = '  +0.38%  '
# Row 28
$ws.Range("E28").Value = '  +0.54%  '
# Row 29
$ws.Range("E29").Value = '  +1.54%  '
# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '0.0464'
$ws.Range("E30").Value = '  +0.68%  '
# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '1.10'
$ws.Range("E31").Value = '  -0.61%  '
# Row 32
$ws.Range("E32").Value = '  +1.18%  '
# Row 33
$ws.Range("D33").Value = '1.405.78'
$ws.Range("E33").Value = '  +4.41%  '
# Row 34
$ws.Range("E34").Value = '  +3.22%  '
# Row 35
$ws.Range("E35").Value = '  +3.70%  '
# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.959'
$ws.Range("E36").Value = '  +3.62%  '
# Row 37
$ws.Range("E37").Value = '  +0.32%  '
# Row 38
$ws.Range("E38").Value = '  +1.66%  '
# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.524'
$ws.Range("E39").Value = '  +0.98%  '
# Row 40
$ws.Range("E40").Value = '  +1.33%  '
# Row 41
$ws.Range("E41").Value = '  +0.54%  '
# Row 42
$ws.Range("E42").Value = '  -0.24%  '
# Row 43
$ws.Range("E43").Value = '  +3.75%  '
# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '5.48'
$ws.Range("E44").Value = '  -3.82%  '
# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '64.02'
$ws.Range("E45").Value = '  +1.94%  '
# Row 46
$ws.Range("E46").Value = '  -0.07%  '
# Row 47
$ws.Range("D47").Value = '1.692.97'
$ws.Range("E47").Value = '  +0.62%  '
# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '86.33'
$ws.Range("E48").Value = '  +0.59%  '
# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.0512'
$ws.Range("E49").Value = '  -1.58%  '
# Row 51
$ws.Range("B51").Value = 'USDD'
$ws.Range("C51").Value = 'https://coinranking.com/coin/z2PZIKQL7+usdd-usdd'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '1.01'
$ws.Range("E51").Value = '  +0.57%  '
